$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 45: add missing formula to G45 (cell already existed but blank)
$ws.Range("G45").Formula = '=B45*$K$6'

# Row 46: add missing formula to G46 (brand new cell - match style of column G)
$ws.Range("G46").Formula = '=B46*$K$6'
$ws.Range("G46").HorizontalAlignment = -4108

# Row 47: add missing formula to G47 (brand new cell - match style of column G)
$ws.Range("G47").Formula = '=B47*$K$6'
$ws.Range("G47").HorizontalAlignment = -4108

# Row 48: add missing formula to G48 (brand new cell - match style of column G)
$ws.Range("G48").Formula = '=B48*$K$6'
$ws.Range("G48").HorizontalAlignment = -4108

# Row 49: brand new row of data
$ws.Range("A49").Value = 43508
$ws.Range("B49").Value = 6
$ws.Range("C49").Value = "Styling"
$ws.Range("D49").Value = 2
$ws.Range("E49").Value = "understanding/Styling Home page"
$ws.Range("F49").Value = "Battering my heid against a wall trying to comprehend CSS and how it relates to this"
$ws.Range("G49").Formula = '=B49*$K$6'
$ws.Range("G49").HorizontalAlignment = -4108

# Row 50: brand new row of data
$ws.Range("A50").Value = 43508
$ws.Range("B50").Value = 2
$ws.Range("C50").Value = "Styling"
$ws.Range("D50").Value = 2
$ws.Range("E50").Value = "Styling Tables/Icons"
$ws.Range("F50").Value = "Getting a bit more confident – afraid to do the layout and break it all"
$ws.Range("G50").Formula = '=B50*$K$6'
$ws.Range("G50").HorizontalAlignment = -4108

# Apply date number format to the new A column cells to match the rest of the date column
$ws.Range("A49:A50").NumberFormat = $ws.Range("A48").NumberFormat

# Update the view state (topLeftCell / selection) to match the edited workbook
$ws.Application.ActiveWindow.ScrollRow = 17
$ws.Range("F51").Select()
